$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Scanner" to "Session"
$ws.Name = "Session"

# Force columns A, C, D to Text format so numeric-looking/date-looking
# strings are stored verbatim instead of being auto-converted to numbers/dates
$ws.Range("A2:A39").NumberFormat = "@"
$ws.Range("C2:C39").NumberFormat = "@"
$ws.Range("D2:D39").NumberFormat = "@"

# Update the scan log rows (Student ID / Log Date / Log Time) with the new values
$ws.Cells.Item(2, 1).Value = "221751"
$ws.Cells.Item(2, 3).Value = "01/10/2025"
$ws.Cells.Item(2, 4).Value = "14:05:19"
$ws.Cells.Item(3, 1).Value = "221655"
$ws.Cells.Item(3, 3).Value = "01/10/2025"
$ws.Cells.Item(3, 4).Value = "14:05:34"
$ws.Cells.Item(4, 1).Value = "221605"
$ws.Cells.Item(4, 3).Value = "01/10/2025"
$ws.Cells.Item(4, 4).Value = "14:06:02"
$ws.Cells.Item(5, 1).Value = "221722"
$ws.Cells.Item(5, 3).Value = "01/10/2025"
$ws.Cells.Item(5, 4).Value = "14:06:14"
$ws.Cells.Item(6, 1).Value = "221641"
$ws.Cells.Item(6, 3).Value = "01/10/2025"
$ws.Cells.Item(6, 4).Value = "14:06:28"
$ws.Cells.Item(7, 1).Value = "221676"
$ws.Cells.Item(7, 3).Value = "01/10/2025"
$ws.Cells.Item(7, 4).Value = "14:06:48"
$ws.Cells.Item(8, 1).Value = "221535"
$ws.Cells.Item(8, 3).Value = "01/10/2025"
$ws.Cells.Item(8, 4).Value = "14:07:02"
$ws.Cells.Item(9, 1).Value = "221712"
$ws.Cells.Item(9, 3).Value = "01/10/2025"
$ws.Cells.Item(9, 4).Value = "14:07:13"
$ws.Cells.Item(10, 1).Value = "221596"
$ws.Cells.Item(10, 3).Value = "01/10/2025"
$ws.Cells.Item(10, 4).Value = "14:07:37"
$ws.Cells.Item(11, 1).Value = "221608"
$ws.Cells.Item(11, 3).Value = "01/10/2025"
$ws.Cells.Item(11, 4).Value = "14:07:57"
$ws.Cells.Item(12, 1).Value = "221700"
$ws.Cells.Item(12, 3).Value = "01/10/2025"
$ws.Cells.Item(12, 4).Value = "14:08:04"
$ws.Cells.Item(13, 1).Value = "221552"
$ws.Cells.Item(13, 3).Value = "01/10/2025"
$ws.Cells.Item(13, 4).Value = "14:09:01"
$ws.Cells.Item(14, 1).Value = "221542"
$ws.Cells.Item(14, 3).Value = "01/10/2025"
$ws.Cells.Item(14, 4).Value = "14:09:48"
$ws.Cells.Item(15, 1).Value = "221719"
$ws.Cells.Item(15, 3).Value = "01/10/2025"
$ws.Cells.Item(15, 4).Value = "14:10:57"
$ws.Cells.Item(16, 1).Value = "221683"
$ws.Cells.Item(16, 3).Value = "01/10/2025"
$ws.Cells.Item(16, 4).Value = "14:11:12"
$ws.Cells.Item(17, 1).Value = "221594"
$ws.Cells.Item(17, 3).Value = "01/10/2025"
$ws.Cells.Item(17, 4).Value = "14:11:30"
$ws.Cells.Item(18, 1).Value = "221697"
$ws.Cells.Item(18, 3).Value = "01/10/2025"
$ws.Cells.Item(18, 4).Value = "14:11:44"
$ws.Cells.Item(19, 1).Value = "221595"
$ws.Cells.Item(19, 3).Value = "01/10/2025"
$ws.Cells.Item(19, 4).Value = "14:12:04"
$ws.Cells.Item(20, 1).Value = "221713"
$ws.Cells.Item(20, 3).Value = "01/10/2025"
$ws.Cells.Item(20, 4).Value = "14:12:57"
$ws.Cells.Item(21, 1).Value = "221533"
$ws.Cells.Item(21, 3).Value = "01/10/2025"
$ws.Cells.Item(21, 4).Value = "14:13:14"
$ws.Cells.Item(22, 1).Value = "221633"
$ws.Cells.Item(22, 3).Value = "01/10/2025"
$ws.Cells.Item(22, 4).Value = "14:13:29"
$ws.Cells.Item(23, 1).Value = "221672"
$ws.Cells.Item(23, 3).Value = "01/10/2025"
$ws.Cells.Item(23, 4).Value = "14:13:50"
$ws.Cells.Item(24, 1).Value = "221702"
$ws.Cells.Item(24, 3).Value = "01/10/2025"
$ws.Cells.Item(24, 4).Value = "14:14:04"
$ws.Cells.Item(25, 1).Value = "221650"
$ws.Cells.Item(25, 3).Value = "01/10/2025"
$ws.Cells.Item(25, 4).Value = "14:14:17"
$ws.Cells.Item(26, 1).Value = "221581"
$ws.Cells.Item(26, 3).Value = "01/10/2025"
$ws.Cells.Item(26, 4).Value = "14:15:54"
$ws.Cells.Item(27, 1).Value = "221583"
$ws.Cells.Item(27, 3).Value = "01/10/2025"
$ws.Cells.Item(27, 4).Value = "14:16:08"
$ws.Cells.Item(28, 1).Value = "221624"
$ws.Cells.Item(28, 3).Value = "01/10/2025"
$ws.Cells.Item(28, 4).Value = "14:16:29"
$ws.Cells.Item(29, 1).Value = "221688"
$ws.Cells.Item(29, 3).Value = "01/10/2025"
$ws.Cells.Item(29, 4).Value = "14:16:40"
$ws.Cells.Item(30, 1).Value = "221686"
$ws.Cells.Item(30, 3).Value = "01/10/2025"
$ws.Cells.Item(30, 4).Value = "14:17:05"
$ws.Cells.Item(31, 1).Value = "221716"
$ws.Cells.Item(31, 3).Value = "01/10/2025"
$ws.Cells.Item(31, 4).Value = "14:17:31"
$ws.Cells.Item(32, 1).Value = "221568"
$ws.Cells.Item(32, 3).Value = "01/10/2025"
$ws.Cells.Item(32, 4).Value = "14:17:46"
$ws.Cells.Item(33, 1).Value = "221652"
$ws.Cells.Item(33, 3).Value = "01/10/2025"
$ws.Cells.Item(33, 4).Value = "14:18:03"
$ws.Cells.Item(34, 1).Value = "221701"
$ws.Cells.Item(34, 3).Value = "01/10/2025"
$ws.Cells.Item(34, 4).Value = "14:18:15"
$ws.Cells.Item(35, 1).Value = "221621"
$ws.Cells.Item(35, 3).Value = "01/10/2025"
$ws.Cells.Item(35, 4).Value = "14:18:30"
$ws.Cells.Item(36, 1).Value = "221615"
$ws.Cells.Item(36, 3).Value = "01/10/2025"
$ws.Cells.Item(36, 4).Value = "14:19:12"
$ws.Cells.Item(37, 1).Value = "221592"
$ws.Cells.Item(37, 3).Value = "01/10/2025"
$ws.Cells.Item(37, 4).Value = "14:19:33"
$ws.Cells.Item(38, 1).Value = "221752"
$ws.Cells.Item(38, 3).Value = "01/10/2025"
$ws.Cells.Item(38, 4).Value = "14:20:04"
$ws.Cells.Item(39, 1).Value = "221740"
$ws.Cells.Item(39, 3).Value = "01/10/2025"
$ws.Cells.Item(39, 4).Value = "14:20:17"

# Drop the trailing rows that no longer exist in the refreshed export (old rows 40-48)
$ws.Range("A40:F48").EntireRow.Delete()

